$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename "Universidad de Granada" -> "UGR" (shared string used across many rows Editorial column)
foreach ($r in (63..104) + (106..107)) {
    $ws.Cells.Item($r, 5).Value = "UGR"
}

# 2. Fill in previously-missing author for row 85 (EM E 23)
$ws.Cells.Item(85, 2).Value = "Francisco Javier Lloréns Montes"

# 3. Clarify title for row 111 (EM E 49)
$ws.Cells.Item(111, 1).Value = "El ejecutivo moderno: dirección comercial"

# 4. Append 11 new catalog rows (EM E 53 .. EM E 63)
# row 115: EM E 53
$ws.Cells.Item(115, 1).Value = "Mundalización económica y crisis político jurídica"
$ws.Cells.Item(115, 2).Value = "Nicolás López Calera"
$ws.Cells.Item(115, 4).Value = 1961.0
$ws.Cells.Item(115, 5).Value = "UGR"
$ws.Cells.Item(115, 6).Value = "EMPRESARIALES"
$ws.Cells.Item(115, 8).Value = "EM E 53"

# row 116: EM E 54
$ws.Cells.Item(116, 1).Value = "Los obstáculos técnicos al comercio en la Comunidad económica europea"
$ws.Cells.Item(116, 2).Value = "Manuel López Escudero"
$ws.Cells.Item(116, 4).Value = 1991.0
$ws.Cells.Item(116, 5).Value = "UGR"
$ws.Cells.Item(116, 6).Value = "EMPRESARIALES"
$ws.Cells.Item(116, 8).Value = "EM E 54"

# row 117: EM E 55
$ws.Cells.Item(117, 1).Value = "La parte maldita"
$ws.Cells.Item(117, 2).Value = "Georges Bataille"
$ws.Cells.Item(117, 4).Value = 1987.0
$ws.Cells.Item(117, 5).Value = "Icaria"
$ws.Cells.Item(117, 6).Value = "EMPRESARIALES"
$ws.Cells.Item(117, 8).Value = "EM E 55"

# row 118: EM E 56
$ws.Cells.Item(118, 1).Value = "Análisis contable del equilibrio financiero de la empresa"
$ws.Cells.Item(118, 2).Value = "Lázaro Rodríguez Ariza"
$ws.Cells.Item(118, 4).Value = 1999.0
$ws.Cells.Item(118, 5).Value = "UGR"
$ws.Cells.Item(118, 6).Value = "EMPRESARIALES"
$ws.Cells.Item(118, 8).Value = "EM E 56"

# row 119: EM E 57
$ws.Cells.Item(119, 1).Value = "Cómo mentir con estadisticas"
$ws.Cells.Item(119, 2).Value = "Darriel Huff"
$ws.Cells.Item(119, 4).Value = 1965.0
$ws.Cells.Item(119, 5).Value = "Sagitario"
$ws.Cells.Item(119, 6).Value = "EMPRESARIALES"
$ws.Cells.Item(119, 8).Value = "EM E 57"

# row 120: EM E 58
$ws.Cells.Item(120, 1).Value = "Las formas ocultas de la propaganda"
$ws.Cells.Item(120, 2).Value = "V. Packard"
$ws.Cells.Item(120, 4).Value = 1959.0
$ws.Cells.Item(120, 5).Value = "EM E 46"
$ws.Cells.Item(120, 6).Value = "EMPRESARIALES"
$ws.Cells.Item(120, 8).Value = "EM E 58"

# row 121: EM E 59
$ws.Cells.Item(121, 1).Value = "Introducción a la historia económica"
$ws.Cells.Item(121, 2).Value = "G. D. H. Cole"
$ws.Cells.Item(121, 4).Value = 1963.0
$ws.Cells.Item(121, 5).Value = "Fondo de cultura económica"
$ws.Cells.Item(121, 6).Value = "EMPRESARIALES"
$ws.Cells.Item(121, 8).Value = "EM E 59"

# row 122: EM E 60
$ws.Cells.Item(122, 1).Value = "El ejecutivo moderno: dirección financiera"
$ws.Cells.Item(122, 2).Value = "Máximo Borrel Vidal"
$ws.Cells.Item(122, 4).Value = 1985.0
$ws.Cells.Item(122, 5).Value = "EM E 49"
$ws.Cells.Item(122, 6).Value = "EMPRESARIALES"
$ws.Cells.Item(122, 8).Value = "EM E 60"

# row 123: EM E 61
$ws.Cells.Item(123, 1).Value = "El ejecutivo moderno: Dirección general"
$ws.Cells.Item(123, 2).Value = "Joaquín Bou Gascons"
$ws.Cells.Item(123, 4).Value = 1985.0
$ws.Cells.Item(123, 5).Value = "EM E 49"
$ws.Cells.Item(123, 6).Value = "EMPRESARIALES"
$ws.Cells.Item(123, 8).Value = "EM E 61"

# row 124: EM E 62
$ws.Cells.Item(124, 1).Value = "El ejecutivo moderno: dirección producción"
$ws.Cells.Item(124, 2).Value = "Francisco Javier Mercader del Campo"
$ws.Cells.Item(124, 4).Value = 1985.0
$ws.Cells.Item(124, 5).Value = "EM E 49"
$ws.Cells.Item(124, 6).Value = "EMPRESARIALES"
$ws.Cells.Item(124, 8).Value = "EM E 62"

# row 125: EM E 63
$ws.Cells.Item(125, 1).Value = "La crisis de las cajas rurales españolas"
$ws.Cells.Item(125, 2).Value = "Antonio Martín Mesa"
$ws.Cells.Item(125, 4).Value = 1988.0
$ws.Cells.Item(125, 5).Value = "UGR"
$ws.Cells.Item(125, 6).Value = "EMPRESARIALES"
$ws.Cells.Item(125, 8).Value = "EM E 63"
